{"js": "// Locate the paragraph ending \"...or find an alternative route\" and:\n//   1) append the missing '?' to finish the sentence, and\n//   2) insert the new \"Break the problem apart\" block right after it\n//      (a blank line, the \"2)\" heading, and four explanatory sentences).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Can he get a bigger boat or find an alternative route\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate target paragraph\");\n}\n\n// 1) Finish the sentence with a question mark.\ntarget.insertText(\"?\", Word.InsertLocation.end);\n\n// 2) Insert the new paragraphs (in reverse order so each lands\n// immediately after the original target paragraph, ending up in the\n// correct forward reading order).\nconst newParagraphs = [\n  \"A sub-goal is to get all entities across the river (If there is water in the river or if it is dried up,  which if it is a dried up river bed then he should just put the cat in the boat and the seeds in his position and push the boat across (which would be hilarious to watch)).  \",\n  \"A sub-gaol is to find a combination that prevents these destructive relationships from occurring using the man\\u2019s presents(If he is required). \",\n  \"The entities themselves are constrains because they have destructive relationships. \",\n  \"The boat itself is a constraint because it limits the method of transportation.\",\n  \"2) Break the problem apart,\",\n];\n\nfor (const text of newParagraphs) {\n  target.insertParagraph(text, Word.InsertLocation.after);\n}\n\n// The leading blank paragraph (directly after the \"?\" sentence): insert a\n// placeholder paragraph and then clear its text so it ends up with a\n// genuine, empty text run (matching the existing blank-paragraph style\n// used throughout this document) rather than no text run at all.\nconst blank = target.insertParagraph(\"\\u00A0\", Word.InsertLocation.after);\nblank.getRange().insertText(\"\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# 1) Append the missing question mark to the end of the \"...alternative\n#    route\" sentence.\n$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Find.Execute(\"find an alternative route\", $false, $false, $false, $false, $false, $true, 1, $false, \"find an alternative route?\", 2) | Out-Null\n\n# 2) Locate that same paragraph again (its text now ends with \"route?\")\n#    and insert the new \"Break the problem apart\" block right after it:\n#    a blank line, the heading, and four explanatory sentences.\n$marker = \"Can he get a bigger boat or find an alternative route?\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($marker)) {\n        $target = $p\n        break\n    }\n}\n\n$apostrophe = [char]0x2019\n$newText = \"`r\" + \"`r2) Break the problem apart,\" + `\n    \"`rThe boat itself is a constraint because it limits the method of transportation.\" + `\n    \"`rThe entities themselves are constrains because they have destructive relationships. \" + `\n    \"`rA sub-gaol is to find a combination that prevents these destructive relationships from occurring using the man${apostrophe}s presents(If he is required). \" + `\n    \"`rA sub-goal is to get all entities across the river (If there is water in the river or if it is dried up,  which if it is a dried up river bed then he should just put the cat in the boat and the seeds in his position and push the boat across (which would be hilarious to watch)).  \"\n\n$target.Range.InsertAfter($newText)\n"}
